# Insert a new weekly data row before the current row 55, shifting the
# existing rows (55-76) down to (56-77), and populate the new row with the
# latest week's price data for "Bruselas (repollito)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55; this pushes rows 55:76 down to 56:77
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with this week's data
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 4).Value = 44704
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 100112035
$ws.Cells.Item(55, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 30
$ws.Cells.Item(55, 11).Value = 30000
$ws.Cells.Item(55, 12).Value = 30000
$ws.Cells.Item(55, 13).Value = 30000
$ws.Cells.Item(55, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 3000
$ws.Cells.Item(55, 17).Value = 10
$ws.Cells.Item(55, 18).Value = "Hortaliza"
